$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.048162677408318
$ws.Range("D2").Value = 1.047296038336663
$ws.Range("E2").Value = 1.051833941173634
$ws.Range("F2").Value = 1.052934466715123
$ws.Range("I2").Value = 1.037749770864372
$ws.Range("J2").Value = 1.053207892268829
$ws.Range("K2").Value = 1.05005907320639
$ws.Range("L2").Value = 1.054584346279199
$ws.Range("M2").Value = 1.055681826644692
$ws.Range("N2").Value = 1.054703568280485
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050395805217226
$ws.Range("D3").Value = 1.049021739318044
$ws.Range("E3").Value = 1.054019005597495
$ws.Range("F3").Value = 1.055136648272025
$ws.Range("I3").Value = 1.038358478188883
$ws.Range("J3").Value = 1.055083492515499
$ws.Range("K3").Value = 1.051594221109864
$ws.Range("L3").Value = 1.056578613733169
$ws.Range("M3").Value = 1.057693395199345
$ws.Range("N3").Value = 1.056581832094639
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051833672705359
$ws.Range("D4").Value = 1.050131896316391
$ws.Range("E4").Value = 1.055426192135116
$ws.Range("F4").Value = 1.056554914025677
$ws.Range("I4").Value = 1.038747735437666
$ws.Range("J4").Value = 1.056289968606323
$ws.Range("K4").Value = 1.052580640488853
$ws.Range("L4").Value = 1.057862002672501
$ws.Range("M4").Value = 1.058987985185605
$ws.Range("N4").Value = 1.057790021519896
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.052436488003893
$ws.Range("D5").Value = 1.050597082971776
$ws.Range("E5").Value = 1.056016207577455
$ws.Range("F5").Value = 1.057149588226205
$ws.Range("I5").Value = 1.038910286363356
$ws.Range("J5").Value = 1.056795490019029
$ws.Range("K5").Value = 1.052993701670049
$ws.Range("L5").Value = 1.058399890200666
$ws.Range("M5").Value = 1.059530583341553
$ws.Range("N5").Value = 1.058296260830993
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.052537606889136
$ws.Range("D6").Value = 1.050675101383285
$ws.Range("E6").Value = 1.056115183025136
$ws.Range("F6").Value = 1.057249345969586
$ws.Range("I6").Value = 1.038937515656469
$ws.Range("J6").Value = 1.056880271676016
$ws.Range("K6").Value = 1.053062961713851
$ws.Range("L6").Value = 1.058490108187823
$ws.Range("M6").Value = 1.059621592371296
$ws.Range("N6").Value = 1.058381162887657
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051841734039469
$ws.Range("D7").Value = 1.050138118109239
$ws.Range("E7").Value = 1.055434082052815
$ws.Range("F7").Value = 1.056562866189422
$ws.Range("I7").Value = 1.03874991172832
$ws.Range("J7").Value = 1.056296729968873
$ws.Range("K7").Value = 1.052586166197861
$ws.Range("L7").Value = 1.057869196385116
$ws.Range("M7").Value = 1.058995241834811
$ws.Range("N7").Value = 1.057796792484357
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048918870546403
$ws.Range("D8").Value = 1.047880608203081
$ws.Range("E8").Value = 1.052573803063868
$ws.Range("F8").Value = 1.053680112902807
$ws.Range("I8").Value = 1.037956449951414
$ws.Range("J8").Value = 1.053843263631506
$ws.Range("K8").Value = 1.050579334319126
$ws.Range("L8").Value = 1.055259796342001
$ws.Range("M8").Value = 1.056363122589712
$ws.Range("N8").Value = 1.055339841943375
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043712131597918
$ws.Range("D9").Value = 1.043851558748543
$ws.Range("E9").Value = 1.047480646105305
$ws.Range("F9").Value = 1.04854736640516
$ws.Range("I9").Value = 1.036522338203573
$ws.Range("J9").Value = 1.049463557497999
$ws.Range("K9").Value = 1.04698873467803
$ws.Range("L9").Value = 1.050606221813939
$ws.Range("M9").Value = 1.051669549093756
$ws.Range("N9").Value = 1.050953916124799
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040200612033816
$ws.Range("D10").Value = 1.041129312864461
$ws.Range("E10").Value = 1.044047192262996
$ws.Range("F10").Value = 1.045087511316246
$ws.Range("I10").Value = 1.035541299771959
$ws.Range("J10").Value = 1.046503683852774
$ws.Range("K10").Value = 1.044556684557471
$ws.Range("L10").Value = 1.047464283421983
$ws.Range("M10").Value = 1.048500952424887
$ws.Range("N10").Value = 1.047989839119496
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038669943602546
$ws.Range("D11").Value = 1.039941531094927
$ws.Range("E11").Value = 1.042550911510465
$ws.Range("F11").Value = 1.043579792791199
$ws.Range("I11").Value = 1.035110395824463
$ws.Range("J11").Value = 1.045212034413165
$ws.Range("K11").Value = 1.04349408397881
$ws.Range("L11").Value = 1.046093898434871
$ws.Range("M11").Value = 1.047119022379989
$ws.Range("N11").Value = 1.046696355389527
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038099809692463
$ws.Range("D12").Value = 1.039498941746466
$ws.Range("E12").Value = 1.041993641018744
$ws.Range("F12").Value = 1.043018272595015
$ws.Range("I12").Value = 1.034949405348841
$ws.Range("J12").Value = 1.044730712773948
$ws.Range("K12").Value = 1.043097922988802
$ws.Range("L12").Value = 1.045583344027716
$ws.Range("M12").Value = 1.046604179040716
$ws.Range("N12").Value = 1.046214350218375
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038222177452442
$ws.Range("D13").Value = 1.039593942369996
$ws.Range("E13").Value = 1.042113245405352
$ws.Range("F13").Value = 1.043138788617549
$ws.Range("I13").Value = 1.034983980857692
$ws.Range("J13").Value = 1.044834028529303
$ws.Range("K13").Value = 1.04318296768592
$ws.Range("L13").Value = 1.045692929747547
$ws.Range("M13").Value = 1.04671468478242
$ws.Range("N13").Value = 1.046317812693953
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038622848585356
$ws.Range("D14").Value = 1.039904975184353
$ws.Range("E14").Value = 1.042504877930882
$ws.Range("F14").Value = 1.043533407956234
$ws.Range("I14").Value = 1.035097107439304
$ws.Range("J14").Value = 1.045172279993289
$ws.Range("K14").Value = 1.043461367257523
$ws.Range("L14").Value = 1.046051727368504
$ws.Range("M14").Value = 1.047076496802844
$ws.Range("N14").Value = 1.046656544513816
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038869505028041
$ws.Range("D15").Value = 1.040096426873497
$ws.Range("E15").Value = 1.042745977429696
$ws.Range("F15").Value = 1.043776347559503
$ws.Range("I15").Value = 1.035166684271638
$ws.Range("J15").Value = 1.045380481802659
$ws.Range("K15").Value = 1.043632703376337
$ws.Range("L15").Value = 1.046272590025694
$ws.Range("M15").Value = 1.047299216644476
$ws.Range("N15").Value = 1.046865041993636
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040301979252586
$ws.Range("D16").Value = 1.04120794856853
$ws.Range("E16").Value = 1.044146289812526
$ws.Range("F16").Value = 1.045187367801891
$ws.Range("I16").Value = 1.035569767410351
$ws.Range("J16").Value = 1.046589191982499
$ws.Range("K16").Value = 1.044627002543867
$ws.Range("L16").Value = 1.047555018913301
$ws.Range("M16").Value = 1.048592454058715
$ws.Range("N16").Value = 1.048075468680572
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041197778465742
$ws.Range("D17").Value = 1.041902732606103
$ws.Range("E17").Value = 1.045022073014033
$ws.Range("F17").Value = 1.046069865963417
$ws.Range("I17").Value = 1.035820964428044
$ws.Range("J17").Value = 1.047344675722442
$ws.Range("K17").Value = 1.045248130034493
$ws.Range("L17").Value = 1.048356769840479
$ws.Range("M17").Value = 1.049400984179667
$ws.Range("N17").Value = 1.048832025294081
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041719306350475
$ws.Range("D18").Value = 1.042307119554178
$ws.Range("E18").Value = 1.045531982145069
$ws.Range("F18").Value = 1.046583691274007
$ws.Range("I18").Value = 1.035966895173436
$ws.Range("J18").Value = 1.047784374214046
$ws.Range("K18").Value = 1.045609508355125
$ws.Range("L18").Value = 1.048823464945287
$ws.Range("M18").Value = 1.049871633316059
$ws.Range("N18").Value = 1.049272348207977
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041896969709689
$ws.Range("D19").Value = 1.042444858807708
$ws.Range("E19").Value = 1.045705693276368
$ws.Range("F19").Value = 1.046758737639043
$ws.Range("I19").Value = 1.036016554490521
$ws.Range("J19").Value = 1.047934138070567
$ws.Range("K19").Value = 1.045732574916264
$ws.Range("L19").Value = 1.048982435657996
$ws.Range("M19").Value = 1.050031952233859
$ws.Range("N19").Value = 1.049422324746352
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04110176906163
$ws.Range("D20").Value = 1.04182827898989
$ws.Range("E20").Value = 1.044928205303
$ws.Range("F20").Value = 1.045975277859237
$ws.Range("I20").Value = 1.035794074322585
$ws.Range("J20").Value = 1.047263719316374
$ws.Range("K20").Value = 1.045181583813845
$ws.Range("L20").Value = 1.048270848408383
$ws.Range("M20").Value = 1.049314335422723
$ws.Range("N20").Value = 1.048750953920631
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038504904816864
$ws.Range("D21").Value = 1.039813422567099
$ws.Range("E21").Value = 1.042389593294188
$ws.Range("F21").Value = 1.04341724381346
$ws.Range("I21").Value = 1.035063820381042
$ws.Range("J21").Value = 1.04507271637689
$ws.Range("K21").Value = 1.043379426220049
$ws.Range("L21").Value = 1.045946113073881
$ws.Range("M21").Value = 1.046969994861299
$ws.Range("N21").Value = 1.046556839505661
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036863012619915
$ws.Range("D22").Value = 1.038538516212073
$ws.Range("E22").Value = 1.04078485031072
$ws.Range("F22").Value = 1.041800282279183
$ws.Range("I22").Value = 1.0345992713124
$ws.Range("J22").Value = 1.043686182300584
$ws.Range("K22").Value = 1.042237850880062
$ws.Range("L22").Value = 1.044475571097128
$ws.Range("M22").Value = 1.045487122769059
$ws.Range("N22").Value = 1.04516833639193
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037734293936083
$ws.Range("D23").Value = 1.039215147366311
$ws.Range("E23").Value = 1.041636387649584
$ws.Range("F23").Value = 1.042658297690473
$ws.Range("I23").Value = 1.034846055717477
$ws.Range("J23").Value = 1.044422074281862
$ws.Range("K23").Value = 1.04284383811513
$ws.Range("L23").Value = 1.045255990779269
$ws.Range("M23").Value = 1.046274079316698
$ws.Range("N23").Value = 1.045905273424228
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041145154594262
$ws.Range("D24").Value = 1.041861924057285
$ws.Range("E24").Value = 1.044970622923059
$ws.Range("F24").Value = 1.046018020995612
$ws.Range("I24").Value = 1.03580622662069
$ws.Range("J24").Value = 1.047300303003444
$ws.Range("K24").Value = 1.045211656011778
$ws.Range("L24").Value = 1.048309675546806
$ws.Range("M24").Value = 1.049353491208411
$ws.Range("N24").Value = 1.048787589560733
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045065124534762
$ws.Range("D25").Value = 1.044899405591427
$ws.Range("E25").Value = 1.048803871511783
$ws.Range("F25").Value = 1.049880826846958
$ws.Range("I25").Value = 1.036897430472652
$ws.Range("J25").Value = 1.050602718258916
$ws.Range("K25").Value = 1.047923609963072
$ws.Range("L25").Value = 1.051816086395515
$ws.Range("M25").Value = 1.052889751829365
$ws.Range("N25").Value = 1.052094694624659
